$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A6").Value = 6
$ws.Range("B6").Value = " five"
$ws.Range("C6").Value = 4
$ws.Range("D6").Value = 53.95540475036354
